# Insert a new weekly record for Femacal de La Calera - Poroto verde.
# This shifts the existing rows 233..241 down to 234..242 and fills the
# newly opened row 233 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(233).Insert()

$ws.Cells.Item(233, 1).Value = 3
$ws.Cells.Item(233, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(233, 3).Value = 'Coquimbo'
$ws.Cells.Item(233, 4).Value = 44509
$ws.Cells.Item(233, 5).Value = 5
$ws.Cells.Item(233, 6).Value = 100112031
$ws.Cells.Item(233, 7).Value = 'Poroto verde'
$ws.Cells.Item(233, 8).Value = 'Magnum'
$ws.Cells.Item(233, 9).Value = 'Primera'
$ws.Cells.Item(233, 10).Value = 38
$ws.Cells.Item(233, 11).Value = 38000
$ws.Cells.Item(233, 12).Value = 38000
$ws.Cells.Item(233, 13).Value = 38000
$ws.Cells.Item(233, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(233, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(233, 16).Value = 1520
$ws.Cells.Item(233, 17).Value = 25
$ws.Cells.Item(233, 18).Value = 'Hortaliza'
